# "Changes to avoid bot detection and fake retailer"
#
# Re-enable the previously-disabled TC001 automation steps that search
# Amazon, type the product name, press Enter, scroll to "Add to cart" and
# click it (rows 10-14, the "Enabled" column) by flipping their value from
# "no" to "Yes". Once "no" is no longer referenced anywhere in the sheet it
# drops out of the shared-strings table automatically on save.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B10:B14").Value = "Yes"

# Clear the stray cell selection that had been left on D7.
$ws.Range("A1").Select()
